$wb = $excel.ActiveWorkbook

# Remove the "genera" and "genera_individual" sheets entirely.
$wb.Worksheets.Item("genera").Delete() | Out-Null
$wb.Worksheets.Item("genera_individual").Delete() | Out-Null

# Update a handful of standard-error estimates on the species_individual sheet
# per reviewer feedback (Nat Commn revision).
$ws = $wb.Worksheets.Item("species_individual")
$ws.Range("E8").Value = 0.5
$ws.Range("E9").Value = 0.38
$ws.Range("E10").Value = 0.78
$ws.Range("E20").Value = 0.61
$ws.Range("E21").Value = 0.46
$ws.Range("E22").Value = 1.12
